$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1190
$ws.Range("I40").Value = 700
$ws.Range("J40").Value = 2333.3333
$ws.Range("K40").Value = 700
$ws.Range("L40").Value = 2333.3333
$ws.Range("M40").Value = -525
$ws.Range("N40").Value = -2683.3333

$ws.Range("H129").Value = 947.2692
$ws.Range("I129").Value = 3173.5
$ws.Range("J129").Value = 858.22
$ws.Range("K129").Value = 9520.5
$ws.Range("L129").Value = 2574.66
$ws.Range("M129").Value = -4520.5
$ws.Range("N129").Value = -12574.66

$ws.Range("H138").Value = 4861.2173
$ws.Range("J138").Value = 5752.9346
$ws.Range("L138").Value = 17258.8038
$ws.Range("N138").Value = -27538.8038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 448.2353
$ws.Range("I97").Value = 466.875
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 466.875
$ws.Range("L97").Value = 150
$ws.Range("M97").Value = 29.125
$ws.Range("N97").Value = -1142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 398.33334
$ws.Range("I22").Value = 398.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 398.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -225.33334
$ws.Range("N22").ClearContents()

$ws.Range("H105").Value = 2115.7827
$ws.Range("I105").Value = 1878.125
$ws.Range("J105").Value = 2659
$ws.Range("K105").Value = 1878.125
$ws.Range("L105").Value = 2659
$ws.Range("M105").Value = -131.125
$ws.Range("N105").Value = -6153

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2131578.5
$ws.Range("I31").Value = 3335597
$ws.Range("K31").Value = 3335597
$ws.Range("M31").Value = -3335302

$ws.Range("H34").Value = 2131578.5
$ws.Range("I34").Value = 3335597
$ws.Range("K34").Value = 3335597
$ws.Range("M34").Value = -3335395

$ws.Range("H74").Value = 24392.166
$ws.Range("J74").Value = 24392.166
$ws.Range("L74").Value = 24392.166
$ws.Range("N74").Value = -26140.166

$ws.Range("H77").Value = 24392.166
$ws.Range("J77").Value = 24392.166
$ws.Range("L77").Value = 73176.49800000001
$ws.Range("N77").Value = -81912.49800000001

$ws.Range("H107").Value = 1279.8
$ws.Range("I107").Value = 1155.2
$ws.Range("K107").Value = 1155.2
$ws.Range("M107").Value = 764.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 386.7
$ws.Range("I4").Value = 173.4
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 520.2
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = -408.2
$ws.Range("N4").Value = -2024

$ws.Range("H75").Value = 1779.6
$ws.Range("I75").Value = 449
$ws.Range("J75").Value = 2666.6667
$ws.Range("K75").Value = 1347
$ws.Range("L75").Value = 8000.000100000001
$ws.Range("M75").Value = -349
$ws.Range("N75").Value = -9996.000100000001

$ws.Range("H78").Value = 1779.6
$ws.Range("I78").Value = 449
$ws.Range("J78").Value = 2666.6667
$ws.Range("K78").Value = 4041
$ws.Range("L78").Value = 24000.0003
$ws.Range("M78").Value = 951
$ws.Range("N78").Value = -33984.0003

$ws.Range("H87").Value = 9192.182000000001
$ws.Range("I87").Value = 5316.2856
$ws.Range("J87").Value = 15975
$ws.Range("K87").Value = 15948.8568
$ws.Range("L87").Value = 47925
$ws.Range("M87").Value = -14700.8568
$ws.Range("N87").Value = -50421

$ws.Range("H90").Value = 9192.182000000001
$ws.Range("I90").Value = 5316.2856
$ws.Range("J90").Value = 15975
$ws.Range("K90").Value = 47846.5704
$ws.Range("L90").Value = 143775
$ws.Range("M90").Value = -41606.5704
$ws.Range("N90").Value = -156255

$ws.Range("H97").Value = 820
$ws.Range("I97").Value = 395
$ws.Range("J97").Value = 1245
$ws.Range("K97").Value = 1185
$ws.Range("L97").Value = 3735
$ws.Range("M97").Value = -689
$ws.Range("N97").Value = -4727

$ws.Range("H98").Value = 459.75
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 546.3333
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 1638.9999
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -4634.9999

$ws.Range("H103").Value = 1937.3889
$ws.Range("I103").Value = 1274.2727
$ws.Range("J103").Value = 2979.4285
$ws.Range("K103").Value = 3822.8181
$ws.Range("L103").Value = 8938.2855
$ws.Range("M103").Value = -2943.8181
$ws.Range("N103").Value = -10696.2855

$ws.Range("H113").Value = 3031199.8
$ws.Range("I113").Value = 25000496
$ws.Range("J113").Value = 951.9655
$ws.Range("K113").Value = 75001488
$ws.Range("L113").Value = 2855.8965
$ws.Range("M113").Value = -74999318
$ws.Range("N113").Value = -7195.8965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83.91304
$ws.Range("I2").Value = 102.77778
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 102.77778
$ws.Range("L2").Value = 16
$ws.Range("M2").Value = 10.22221999999999
$ws.Range("N2").Value = -242

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H57").Value = 19730.5
$ws.Range("J57").Value = 23845.75
$ws.Range("L57").Value = 23845.75
$ws.Range("N57").Value = -25485.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2572.0208
$ws.Range("I132").Value = 1801.7273
$ws.Range("J132").Value = 4266.6665
$ws.Range("K132").Value = 5405.1819
$ws.Range("L132").Value = 12799.9995
$ws.Range("M132").Value = -2875.1819
$ws.Range("N132").Value = -17859.9995

$ws.Range("H136").Value = 3577719
$ws.Range("I136").Value = 8339701
$ws.Range("J136").Value = 6232.5
$ws.Range("K136").Value = 25019103
$ws.Range("L136").Value = 18697.5
$ws.Range("M136").Value = -25016553
$ws.Range("N136").Value = -23797.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3062.875
$ws.Range("I107").Value = 917.1667
$ws.Range("J107").Value = 9500
$ws.Range("K107").Value = 2751.5001
$ws.Range("L107").Value = 28500
$ws.Range("M107").Value = -831.5001000000002
$ws.Range("N107").Value = -32340
